$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) values.
# D-column values are forced to Text format before assignment so that
# numeric-looking strings (e.g. "312.58") are not reinterpreted by Excel
# as numbers; the style is then reset back to Normal so no stray cell
# style is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.096.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5042"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3832"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08494"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.116"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.251"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.888.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.223"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06660"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.091"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.135.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.264"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.589"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.094.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "156.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("E31").Value = "  -0.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.053"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.638"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.617"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.733"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02451"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06554"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2178"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6517"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.248"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.916"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6192"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.298"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.681"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("E48").Value = "  +1.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.223"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.93%  "
